# Update the "想去人数" (interest count) figures in column F for the two
# sheets that carry the full event listing ("展览" and "全部类型").
# The row numbers correspond directly to the worksheet row numbers.

$wb = $excel.ActiveWorkbook

$updates = @{
    3  = 376
    12 = 1142
    13 = 1468
    14 = 326
    17 = 101
    19 = 57
    30 = 107
    31 = 3981
    32 = 10
    35 = 1017
    36 = 107
    39 = 116
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
